$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches original inlineStr formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "36.383.11"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "2.039.28"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "245.04"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "0.656"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "53.02"
$ws.Range("E8").Value = "  -8.39%  "
$ws.Range("D9").Value = "61.74"
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("D10").Value = "0.358"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").Value = "0.0736"
$ws.Range("E11").Value = "  -5.62%  "
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("D13").Value = "0.916"
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("D14").Value = "14.36"
$ws.Range("E14").Value = "  -5.74%  "
$ws.Range("D15").Value = "2.340.34"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "5.33"
$ws.Range("E16").Value = "  -5.13%  "
$ws.Range("D17").Value = "2.043.01"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "36.245.15"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "16.74"
$ws.Range("E19").Value = "  -6.75%  "
$ws.Range("D20").Value = "70.83"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").Value = "0.0₃0845"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").Value = "234.99"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "5.10"
$ws.Range("E23").Value = "  -5.59%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("D27").Value = "163.26"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("D28").Value = "9.01"
$ws.Range("E28").Value = "  -12.87%  "
$ws.Range("D29").Value = "19.65"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.91"
$ws.Range("E31").Value = "  -10.03%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("D33").Value = "0.0585"
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("D34").Value = "4.34"
$ws.Range("E34").Value = "  -7.61%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "0.0867"
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("D37").Value = "1.81"
$ws.Range("D38").Value = "2.18"
$ws.Range("E38").Value = "  -7.17%  "
$ws.Range("D39").Value = "4.92"
$ws.Range("E39").Value = "  -4.43%  "
$ws.Range("D40").Value = "1.21"
$ws.Range("E40").Value = "  -8.16%  "
$ws.Range("E41").Value = "  -4.38%  "
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("D43").Value = "1.08"
$ws.Range("E43").Value = "  -5.12%  "
$ws.Range("D44").Value = "92.13"
$ws.Range("E44").Value = "  -4.78%  "
$ws.Range("D45").Value = "0.0889"
$ws.Range("E45").Value = "  -6.14%  "
$ws.Range("D46").Value = "1.370.34"
$ws.Range("E46").Value = "  +5.40%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "7.33"
$ws.Range("E47").Value = "  +8.76%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "15.42"
$ws.Range("E48").Value = "  -8.73%  "
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").Value = "2.224.25"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "2.21"
$ws.Range("E51").Value = "  -6.11%  "
